$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Output GWh) values for rows 2-35
$ws.Cells.Item(2, 3).Value = 59.300245142098426
$ws.Cells.Item(3, 3).Value = 472.3506590243318
$ws.Cells.Item(4, 3).Value = 206.27609135887923
$ws.Cells.Item(5, 3).Value = 102.63138930362442
$ws.Cells.Item(6, 3).Value = 10.034425619606871
$ws.Cells.Item(7, 3).Value = 351.1017287873133
$ws.Cells.Item(8, 3).Value = 53.70110086252858
$ws.Cells.Item(9, 3).Value = 151.50694721269838
$ws.Cells.Item(10, 3).Value = 86.70118643976329
$ws.Cells.Item(11, 3).Value = 54.53718464232486
$ws.Cells.Item(12, 3).Value = 8.553909032751054
$ws.Cells.Item(13, 3).Value = 47.135704722896705
$ws.Cells.Item(14, 3).Value = 10.530942327024981
$ws.Cells.Item(15, 3).Value = 204.72481153124102
$ws.Cells.Item(16, 3).Value = 44.21177954911971
$ws.Cells.Item(17, 3).Value = 1122.5543353582025
$ws.Cells.Item(18, 3).Value = 16.72051401459426
$ws.Cells.Item(19, 3).Value = 21.536318999774704
$ws.Cells.Item(20, 3).Value = 427.3246897714147
$ws.Cells.Item(21, 3).Value = 40.68262109504845
$ws.Cells.Item(22, 3).Value = 370.9451144353257
$ws.Cells.Item(23, 3).Value = 70.28142033112843
$ws.Cells.Item(24, 3).Value = 31.57875399031387
$ws.Cells.Item(25, 3).Value = 209.2514383590021
$ws.Cells.Item(26, 3).Value = 30.20114586205253
$ws.Cells.Item(27, 3).Value = 158.7695369417144
$ws.Cells.Item(28, 3).Value = 29.028164074802728
$ws.Cells.Item(29, 3).Value = 522.658652504415
$ws.Cells.Item(30, 3).Value = 119.24508429068537
$ws.Cells.Item(31, 3).Value = 548.4981806502778
$ws.Cells.Item(32, 3).Value = 6.574229050401869
$ws.Cells.Item(33, 3).Value = 35.0061962168337
$ws.Cells.Item(34, 3).Value = 7.5039038344233155
$ws.Cells.Item(35, 3).Value = 5631.658405336615

# Update column E (Operational Capacity GW) values for rows 2-35
$ws.Cells.Item(2, 5).Value = 0.009926797806487257
$ws.Cells.Item(3, 5).Value = 0.24052569297081206
$ws.Cells.Item(4, 5).Value = 0.07444064369848587
$ws.Cells.Item(5, 5).Value = 0.061815767457659786
$ws.Cells.Item(6, 5).Value = 0.004349694603477011
$ws.Cells.Item(7, 5).Value = 0.2812898480369782
$ws.Cells.Item(8, 5).Value = 0.022965901923440005
$ws.Cells.Item(9, 5).Value = 0.08244949855582004
$ws.Cells.Item(10, 5).Value = 0.02015473934116426
$ws.Cells.Item(11, 5).Value = 0.03373750362443276
$ws.Cells.Item(12, 5).Value = 0.005182693663728391
$ws.Cells.Item(13, 5).Value = 0.019736773725196114
$ws.Cells.Item(14, 5).Value = 0.005544676179743804
$ws.Cells.Item(15, 5).Value = 0.2336242468817131
$ws.Cells.Item(16, 5).Value = 0.028469548967120173
$ws.Cells.Item(17, 5).Value = 0.3611119565544601
$ws.Cells.Item(18, 5).Value = 0.012085878771943427
$ws.Cells.Item(19, 5).Value = 0.01350686245337584
$ws.Cells.Item(20, 5).Value = 0.09433601148009638
$ws.Cells.Item(21, 5).Value = 0.021994571400503645
$ws.Cells.Item(22, 5).Value = 0.12388620686361766
$ws.Cells.Item(23, 5).Value = 0.017758489389976276
$ws.Cells.Item(24, 5).Value = 0.016577465964747767
$ws.Cells.Item(25, 5).Value = 0.10092258930706463
$ws.Cells.Item(26, 5).Value = 0.022668711662704042
$ws.Cells.Item(27, 5).Value = 0.07867127891076915
$ws.Cells.Item(28, 5).Value = 0.014793001534613054
$ws.Cells.Item(29, 5).Value = 0.08093261078174376
$ws.Cells.Item(30, 5).Value = 0.06896442298645364
$ws.Cells.Item(31, 5).Value = 0.1072088794755035
$ws.Cells.Item(32, 5).Value = 0.008323954911969252
$ws.Cells.Item(33, 5).Value = 0.019421545766620484
$ws.Cells.Item(34, 5).Value = 0.00417650661018236
$ws.Cells.Item(35, 5).Value = 2.291554972262604

# Column D (Output %) becomes a formula: D2 standalone, D3:D35 shared
$ws.Range("D2").Formula = "=C2/`$C`$35"
$ws.Range("D3:D35").Formula = "=C3/`$C`$35"

# Column F (Capacity %) becomes a formula: F2 standalone, F3:F35 shared
$ws.Range("F2").Formula = "=E2/`$E`$35"
$ws.Range("F3:F35").Formula = "=E3/`$E`$35"

# Update the active selection to C38 (as in the diff)
$ws.Range("C38").Select()
